# Commit: "upload vendor pincode changes"
# The "State" / "{sc:state}" column (column B) is no longer needed in the
# Service Charges list export, so remove it entirely. Excel will shift all
# subsequent columns one position to the left and shrink the used range
# from A1:K2 down to A1:J2.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the entire "State" column (column B); remaining columns (Product,
# Category, Capacity, Service Category, Vendor Basic Charge, Vendor Tax,
# Vendor Total, Customer Total Rs., Serial Number Mandatory, and their
# matching {sc:...} placeholder row) shift left automatically.
$ws.Columns.Item(2).Delete()

# Restore a sensible view: select D8 (mirroring the original author's
# cursor position after removing a column) and make sure column A is back
# in view (the sheet had previously been scrolled right to show column F).
$ws.Range("D8").Select()
$excel.ActiveWindow.ScrollColumn = 1
